$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Formula = "=""" + $value + """"
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range("D2") "26.993.53"
Set-TextValue $ws.Range("E2") "  -2.83%  "
Set-TextValue $ws.Range("D3") "1.741.56"
Set-TextValue $ws.Range("E3") "  -1.03%  "
Set-TextValue $ws.Range("D4") "0.9994"
Set-TextValue $ws.Range("E4") "  -0.38%  "
Set-TextValue $ws.Range("D5") "310.83"
Set-TextValue $ws.Range("E5") "  -5.37%  "
Set-TextValue $ws.Range("D6") "0.9991"
Set-TextValue $ws.Range("E6") "  -0.19%  "
Set-TextValue $ws.Range("D7") "0.4976"
Set-TextValue $ws.Range("E7") "  +3.56%  "
Set-TextValue $ws.Range("D8") "0.3568"
Set-TextValue $ws.Range("E8") "  +0.45%  "
Set-TextValue $ws.Range("D9") "42.58"
Set-TextValue $ws.Range("E9") "  -1.30%  "
Set-TextValue $ws.Range("D10") "0.07264"
Set-TextValue $ws.Range("E10") "  -3.25%  "
Set-TextValue $ws.Range("D11") "1.062"
Set-TextValue $ws.Range("E11") "  -1.62%  "
Set-TextValue $ws.Range("D12") "0.9992"
Set-TextValue $ws.Range("E12") "  -0.17%  "
Set-TextValue $ws.Range("D13") "20.02"
Set-TextValue $ws.Range("E13") "  -2.53%  "
Set-TextValue $ws.Range("D14") "5.982"
Set-TextValue $ws.Range("E14") "  -1.74%  "
Set-TextValue $ws.Range("D15") "1.738.44"
Set-TextValue $ws.Range("E15") "  -1.23%  "
Set-TextValue $ws.Range("D16") "6.856"
Set-TextValue $ws.Range("E16") "  -3.97%  "
Set-TextValue $ws.Range("D17") "86.54"
Set-TextValue $ws.Range("E17") "  -6.83%  "
Set-TextValue $ws.Range("E18") "  -4.80%  "
Set-TextValue $ws.Range("D19") "0.06391"
Set-TextValue $ws.Range("E19") "  -0.63%  "
Set-TextValue $ws.Range("D20") "0.9996"
Set-TextValue $ws.Range("E20") "  -0.12%  "
Set-TextValue $ws.Range("D21") "16.61"
Set-TextValue $ws.Range("E21") "  -1.07%  "
Set-TextValue $ws.Range("D22") "5.744"
Set-TextValue $ws.Range("E22") "  -0.95%  "
Set-TextValue $ws.Range("D23") "27.042.03"
Set-TextValue $ws.Range("E23") "  -2.83%  "
Set-TextValue $ws.Range("D24") "11.23"
Set-TextValue $ws.Range("E24") "  +1.00%  "
Set-TextValue $ws.Range("D25") "2.053"
Set-TextValue $ws.Range("E25") "  -4.96%  "
Set-TextValue $ws.Range("D26") "153.45"
Set-TextValue $ws.Range("E26") "  -6.14%  "
Set-TextValue $ws.Range("D27") "19.90"
Set-TextValue $ws.Range("E27") "  -0.99%  "
Set-TextValue $ws.Range("D28") "1.933.81"
Set-TextValue $ws.Range("E28") "  -1.51%  "
Set-TextValue $ws.Range("E29") "  -3.37%  "
Set-TextValue $ws.Range("D30") "120.93"
Set-TextValue $ws.Range("E30") "  -1.68%  "
Set-TextValue $ws.Range("D31") "1.063"
Set-TextValue $ws.Range("E31") "  +0.75%  "
Set-TextValue $ws.Range("D32") "0.09444"
Set-TextValue $ws.Range("E32") "  -0.09%  "
Set-TextValue $ws.Range("D33") "3.570"
Set-TextValue $ws.Range("E33") "  -2.23%  "
Set-TextValue $ws.Range("D34") "5.385"
Set-TextValue $ws.Range("E34") "  -2.82%  "
Set-TextValue $ws.Range("D35") "0.02205"
Set-TextValue $ws.Range("E35") "  -2.70%  "
Set-TextValue $ws.Range("D36") "0.05919"
Set-TextValue $ws.Range("E36") "  -0.95%  "
Set-TextValue $ws.Range("D37") "11.09"
Set-TextValue $ws.Range("E37") "  -4.50%  "
Set-TextValue $ws.Range("B38") "WEMIXTOKEN"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D38") "1.422"
Set-TextValue $ws.Range("E38") "  -0.94%  "
Set-TextValue $ws.Range("B39") "Algorand"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D39") "0.1997"
Set-TextValue $ws.Range("E39") "  -3.24%  "
Set-TextValue $ws.Range("D40") "4.766"
Set-TextValue $ws.Range("E40") "  -2.34%  "
Set-TextValue $ws.Range("B41") "TheSandbox"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D41") "0.6016"
Set-TextValue $ws.Range("E41") "  -1.99%  "
Set-TextValue $ws.Range("B42") "Frax"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D42") "0.9993"
Set-TextValue $ws.Range("E42") "  -0.11%  "
Set-TextValue $ws.Range("D43") "1.112"
Set-TextValue $ws.Range("E43") "  -6.08%  "
Set-TextValue $ws.Range("D44") "7.478"
Set-TextValue $ws.Range("E44") "  -3.69%  "
Set-TextValue $ws.Range("D45") "12.83"
Set-TextValue $ws.Range("E45") "  -2.48%  "
Set-TextValue $ws.Range("D46") "3.581"
Set-TextValue $ws.Range("E46") "  -4.10%  "
Set-TextValue $ws.Range("D47") "0.5647"
Set-TextValue $ws.Range("E47") "  -2.54%  "
Set-TextValue $ws.Range("D48") "120.05"
Set-TextValue $ws.Range("E48") "  -2.51%  "
Set-TextValue $ws.Range("D49") "1.861"
Set-TextValue $ws.Range("E49") "  -3.05%  "
Set-TextValue $ws.Range("D50") "1.102"
Set-TextValue $ws.Range("E50") "  -3.73%  "
Set-TextValue $ws.Range("D51") "0.06667"
Set-TextValue $ws.Range("E51") "  -1.71%  "
